$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the new date values, and apply a date/time
# number format to column G (style carried over from formatting the column).
$ws.Columns.Item(7).NumberFormat = "m/d/yy h:mm"

# New trade row
$ws.Range("A3").Value = 42636.592731481483
$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
$ws.Range("B3").Value = $true
$ws.Range("C3").Value = 10030.5
$ws.Range("D3").Value = 10000
$ws.Range("E3").Value = 81.97
$ws.Range("F3").Value = 81.47
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = -0.61
$ws.Range("I3").Value = $false

$ws.Columns.Item(1).AutoFit()
